$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update utilisation cost values in column F
$ws.Range("F2").Value = 858.82
$ws.Range("F3").Value = 378.83
$ws.Range("F4").Value = 128.1
$ws.Range("F6").Value = 751.9
$ws.Range("F7").Value = 751.9
$ws.Range("F8").Value = 751.9
$ws.Range("F9").Value = 1938.76

# Update selection to F10 (matches the author's final cursor position)
$ws.Range("F10").Select()
